# Update cryptos list (price & volume figures) to match latest snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.300.13'
$ws.Range("E2").Value = '  +2.76%  '

$ws.Range("D3").Value = '2.264.46'
$ws.Range("E3").Value = '  +1.65%  '

$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '320.23'
$ws.Range("E5").Value = '  -0.58%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '102.57'
$ws.Range("E6").Value = '  +3.67%  '

$ws.Range("E7").Value = '  -0.19%  '

$ws.Range("E8").Value = '  +0.10%  '

$ws.Range("E9").Value = '  -1.67%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.39'
$ws.Range("E10").Value = '  +1.24%  '

$ws.Range("E11").Value = '  +0.89%  '

$ws.Range("E12").Value = '  +0.12%  '

$ws.Range("E13").Value = '  -0.96%  '

$ws.Range("D14").Value = '2.599.69'
$ws.Range("E14").Value = '  +1.32%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.864'
$ws.Range("E15").Value = '  +0.08%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.52'

$ws.Range("D17").Value = '2.262.64'
$ws.Range("E17").Value = '  +1.96%  '

$ws.Range("D18").Value = '44.171.78'
$ws.Range("E18").Value = '  +2.71%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.52'
$ws.Range("E19").Value = '  -3.09%  '

$ws.Range("E20").Value = '  +2.36%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.55'
$ws.Range("E21").Value = '  +0.04%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.91'
$ws.Range("E22").Value = '  +0.89%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.15'
$ws.Range("E23").Value = '  -2.63%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '236.45'
$ws.Range("E24").Value = '  -0.43%  '

$ws.Range("E25").Value = '  -4.03%  '

$ws.Range("E26").Value = '  +0.17%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.59'
$ws.Range("E27").Value = '  +5.65%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '38.93'
$ws.Range("E28").Value = '  +7.17%  '

$ws.Range("E29").Value = '  -1.68%  '

$ws.Range("E30").Value = '  -2.20%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '163.49'
$ws.Range("E31").Value = '  +5.45%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.32'
$ws.Range("E32").Value = '  -0.21%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0854'
$ws.Range("E33").Value = '  -1.88%  '

$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.01'
$ws.Range("E35").Value = '  +5.54%  '

$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.114'
$ws.Range("E36").Value = '  +9.32%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.06'
$ws.Range("E37").Value = '  -7.37%  '

$ws.Range("E38").Value = '  -1.40%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '16.63'
$ws.Range("E39").Value = '  +17.19%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.74'
$ws.Range("E40").Value = '  +0.92%  '

$ws.Range("E41").Value = '  -4.24%  '

$ws.Range("E42").Value = '  -1.54%  '

$ws.Range("E43").Value = '  +0.12%  '

$ws.Range("D44").Value = '1.782.06'
$ws.Range("E44").Value = '  +2.54%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.199'
$ws.Range("E45").Value = '  -1.85%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '82.79'
$ws.Range("E46").Value = '  -2.69%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '75.54'
$ws.Range("E47").Value = '  +0.38%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.25'
$ws.Range("E48").Value = '  -0.64%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '105.10'
$ws.Range("E49").Value = '  +1.92%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '58.58'
$ws.Range("E50").Value = '  +1.04%  '

$ws.Range("E51").Value = '  +4.62%  '
